# "Sửa thông tin cá nhân" – a pending "Ban Đào Tạo" submission (ID 123345,
# previously dated 17/09/2022, status "Chưa duyệt") is corrected and sent
# up for approval: its date becomes 18/09/2022, its status becomes
# "Chờ duyệt", and it now carries an explicit SubID of "N/A". The same
# (corrected) record is appended as a brand-new row to the
# "Ban Điều Hành duyệt" worksheet so the approving board can review it.

$wb = $excel.ActiveWorkbook

$wsDuyet = $wb.Worksheets.Item(1)   # "Ban Điều Hành duyệt"
$wsDaoTao = $wb.Worksheets.Item(5)  # "Ban Đào Tạo"

# --- 1. Update the source record on "Ban Đào Tạo" -------------------------
$wsDaoTao.Range("B3").Value2 = "N/A"
$wsDaoTao.Range("E3").Value2 = "18/09/2022"
$wsDaoTao.Range("N3").Value2 = "Chờ duyệt"

# the pending-change event log gains a new "hang xoa = 2" entry, pushing
# the existing queue down by one row
$wsDaoTao.Range("J9").Value2 = "hang xoa = 3"
$wsDaoTao.Range("J8").Value2 = "hang xoa = 2"

# --- 2. Renumber the existing rows on "Ban Điều Hành duyệt" ---------------
$wsDuyet.Range("A2").Value2 = 5
$wsDuyet.Range("A3").Value2 = 4
$wsDuyet.Range("A4").Value2 = 3
$wsDuyet.Range("A5").Value2 = 2

# --- 3. Append the new row (the corrected submission) ----------------------
$wsDuyet.Range("A6").Value2 = 1
$wsDuyet.Range("B6").Value2 = "N/A"
$wsDuyet.Range("C6").Value2 = "123345"
$wsDuyet.Range("D6").Value2 = "Ban Đào Tạo"
$wsDuyet.Range("E6").Value2 = "18/09/2022"
$wsDuyet.Range("F6").Value2 = "<p>1235</p>"
$wsDuyet.Range("G6").Value2 = "https://www.plus2net.com"
$wsDuyet.Range("H6").Value2 = $false
$wsDuyet.Range("I6").Value2 = "Chưa có phản hồi"
$wsDuyet.Range("J6").Value2 = $false
$wsDuyet.Range("K6").Value2 = "Chưa có phản hồi"
$wsDuyet.Range("L6").Value2 = $false
$wsDuyet.Range("N6").Value2 = "Chờ duyệt"

# --- 4. The pending-change event log (column J, below the data rows) shifts
#        up by one: the "subid=BNS4" entry that used to sit at J6 is
#        consumed by the new row above, a new "hang xoa = 4" entry is
#        recorded at J7, and the remaining entries shift up, finally
#        dropping the trailing "subid=N/A" row altogether.
$wsDuyet.Range("J10").Value2 = "subid=N/A"
$wsDuyet.Range("J9").Value2 = "hang xoa = 6"
$wsDuyet.Range("J8").Value2 = "hang xoa = 6"
$wsDuyet.Range("J7").Value2 = "hang xoa = 4"
$wsDuyet.Range("J11").ClearContents()
